$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 458.2
$ws.Range("J17").Value = 458.2
$ws.Range("L17").Value = 1374.6
$ws.Range("N17").Value = -1710.6
$ws.Range("H19").Value = 651.2
$ws.Range("I19").Value = 536
$ws.Range("K19").Value = 536
$ws.Range("M19").Value = -361
$ws.Range("H38").Value = 287.88235
$ws.Range("I38").Value = 62.125
$ws.Range("K38").Value = 186.375
$ws.Range("M38").Value = 185.625
$ws.Range("H39").Value = 275.73914
$ws.Range("I39").Value = 207.13333
$ws.Range("J39").Value = 404.375
$ws.Range("K39").Value = 621.39999
$ws.Range("L39").Value = 1213.125
$ws.Range("M39").Value = -325.39999
$ws.Range("N39").Value = -1805.125
$ws.Range("H70").Value = 23953638
$ws.Range("I70").Value = 41917744
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 125753232
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -125752962
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 23953638
$ws.Range("I73").Value = 41917744
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 125753232
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -125752296
$ws.Range("N73").Value = -6372
$ws.Range("H86").Value = 3872.5
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 4396
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 4396
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -6642
$ws.Range("H88").Value = 1620.25
$ws.Range("J88").Value = 2263.25
$ws.Range("L88").Value = 2263.25
$ws.Range("N88").Value = -3075.25
$ws.Range("H89").Value = 3872.5
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 4396
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 21980
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -33212
$ws.Range("H91").Value = 1620.25
$ws.Range("J91").Value = 2263.25
$ws.Range("L91").Value = 2263.25
$ws.Range("N91").Value = -5071.25
$ws.Range("H112").Value = 2972.7334
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 3113.6428
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 9340.928400000001
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -11556.9284
$ws.Range("H129").Value = 1068.9066
$ws.Range("J129").Value = 1212
$ws.Range("L129").Value = 3636
$ws.Range("N129").Value = -13636
$ws.Range("H138").Value = 3585.8667
$ws.Range("I138").Value = 1591.6522
$ws.Range("J138").Value = 4825.5137
$ws.Range("K138").Value = 4774.9566
$ws.Range("L138").Value = 14476.5411
$ws.Range("M138").Value = 365.0434000000005
$ws.Range("N138").Value = -24756.5411
$ws.Range("H141").Value = 7476.8096
$ws.Range("I141").Value = 2717.111
$ws.Range("K141").Value = 8151.333
$ws.Range("M141").Value = -2971.333

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2163.3845
$ws.Range("I2").Value = 1679
$ws.Range("K2").Value = 1679
$ws.Range("M2").Value = -1566
$ws.Range("H74").Value = 2467.889
$ws.Range("I74").Value = 2509.5217
$ws.Range("J74").Value = 2228.5
$ws.Range("K74").Value = 2509.5217
$ws.Range("L74").Value = 2228.5
$ws.Range("M74").Value = -1635.5217
$ws.Range("N74").Value = -3976.5
$ws.Range("H77").Value = 2467.889
$ws.Range("I77").Value = 2509.5217
$ws.Range("J77").Value = 2228.5
$ws.Range("K77").Value = 12547.6085
$ws.Range("L77").Value = 11142.5
$ws.Range("M77").Value = -8179.608499999998
$ws.Range("N77").Value = -19878.5
$ws.Range("H113").Value = 37232.668
$ws.Range("J113").Value = 37232.668
$ws.Range("L113").Value = 37232.668
$ws.Range("N113").Value = -45910.668
$ws.Range("H116").Value = 2163.3845
$ws.Range("I116").Value = 1679
$ws.Range("K116").Value = 1679
$ws.Range("M116").Value = 615

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2163.3845
$ws.Range("I3").Value = 1679
$ws.Range("K3").Value = 1679
$ws.Range("M3").Value = -1565
$ws.Range("H20").Value = 112275.11
$ws.Range("I20").Value = 167828
$ws.Range("J20").Value = 1169.3334
$ws.Range("K20").Value = 167828
$ws.Range("L20").Value = 1169.3334
$ws.Range("M20").Value = -167581
$ws.Range("N20").Value = -1663.3334
$ws.Range("H94").Value = 794.0645
$ws.Range("I94").Value = 751.5
$ws.Range("J94").Value = 940
$ws.Range("K94").Value = 751.5
$ws.Range("L94").Value = 940
$ws.Range("M94").Value = -300.5
$ws.Range("N94").Value = -1842

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4183.7383
$ws.Range("I31").Value = 1839.6875
$ws.Range("J31").Value = 4949.143
$ws.Range("K31").Value = 1839.6875
$ws.Range("L31").Value = 4949.143
$ws.Range("M31").Value = -1544.6875
$ws.Range("N31").Value = -5539.143
$ws.Range("H34").Value = 4183.7383
$ws.Range("I34").Value = 1839.6875
$ws.Range("J34").Value = 4949.143
$ws.Range("K34").Value = 1839.6875
$ws.Range("L34").Value = 4949.143
$ws.Range("M34").Value = -1637.6875
$ws.Range("N34").Value = -5353.143

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1233.8462
$ws.Range("J33").Value = 1428.8889
$ws.Range("L33").Value = 8573.3334
$ws.Range("N33").Value = -9139.3334
$ws.Range("H64").Value = 3142.8667
$ws.Range("J64").Value = 4234.222
$ws.Range("L64").Value = 12702.666
$ws.Range("N64").Value = -13242.666
$ws.Range("H67").Value = 3142.8667
$ws.Range("J67").Value = 4234.222
$ws.Range("L67").Value = 12702.666
$ws.Range("N67").Value = -14574.666
$ws.Range("H131").Value = 34875.184
$ws.Range("I131").Value = 756
$ws.Range("J131").Value = 42629.547
$ws.Range("K131").Value = 2268
$ws.Range("L131").Value = 127888.641
$ws.Range("M131").Value = 2772
$ws.Range("N131").Value = -137968.641
$ws.Range("H133").Value = 4074.238
$ws.Range("I133").Value = 1868.1
$ws.Range("J133").Value = 6079.8184
$ws.Range("K133").Value = 5604.299999999999
$ws.Range("L133").Value = 18239.4552
$ws.Range("M133").Value = -544.2999999999993
$ws.Range("N133").Value = -28359.4552

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2269
$ws.Range("I61").Value = 2116.4707
$ws.Range("K61").Value = 2116.4707
$ws.Range("M61").Value = -1914.4707
$ws.Range("H113").Value = 2269
$ws.Range("I113").Value = 2116.4707
$ws.Range("K113").Value = 2116.4707
$ws.Range("M113").Value = 53.52930000000015
$ws.Range("H132").Value = 5753.2
$ws.Range("I132").Value = 6755.5557
$ws.Range("J132").Value = 4249.6665
$ws.Range("K132").Value = 20266.6671
$ws.Range("L132").Value = 12748.9995
$ws.Range("M132").Value = -17736.6671
$ws.Range("N132").Value = -17808.9995

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 42347.5
$ws.Range("I57").Value = 52000
$ws.Range("J57").Value = 39130
$ws.Range("K57").Value = 52000
$ws.Range("L57").Value = 39130
$ws.Range("M57").Value = -51246
$ws.Range("N57").Value = -40638
$ws.Range("H62").Value = 3950
$ws.Range("I62").Value = 3950
$ws.Range("K62").Value = 3950
$ws.Range("M62").Value = -3326
$ws.Range("H65").Value = 3950
$ws.Range("I65").Value = 3950
$ws.Range("K65").Value = 19750
$ws.Range("M65").Value = -16630
$ws.Range("H81").Value = 59373.15
$ws.Range("J81").Value = 4353.8184
$ws.Range("L81").Value = 8707.6368
$ws.Range("N81").Value = -10829.6368
$ws.Range("H84").Value = 59373.15
$ws.Range("J84").Value = 4353.8184
$ws.Range("L84").Value = 43538.184
$ws.Range("N84").Value = -54146.184
$ws.Range("H87").Value = 84618.75
$ws.Range("J87").Value = 25325
$ws.Range("L87").Value = 25325
$ws.Range("N87").Value = -27821
$ws.Range("H90").Value = 84618.75
$ws.Range("J90").Value = 25325
$ws.Range("L90").Value = 75975
$ws.Range("N90").Value = -88455
